$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 53264.047
$ws.Range("J17").Value = 53264.047
$ws.Range("L17").Value = 159792.141
$ws.Range("N17").Value = -160128.141
$ws.Range("H87").Value = 27780.08
$ws.Range("J87").Value = 27780.08
$ws.Range("L87").Value = 27780.08
$ws.Range("N87").Value = -30276.08
$ws.Range("H90").Value = 27780.08
$ws.Range("J90").Value = 27780.08
$ws.Range("L90").Value = 83340.24000000001
$ws.Range("N90").Value = -95820.24000000001
$ws.Range("H98").Value = 2986
$ws.Range("I98").Value = 3178.111
$ws.Range("J98").Value = 1833.3334
$ws.Range("K98").Value = 3178.111
$ws.Range("L98").Value = 1833.3334
$ws.Range("M98").Value = -1680.111
$ws.Range("N98").Value = -4829.3334
$ws.Range("H122").Value = 2986
$ws.Range("I122").Value = 3178.111
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 9534.332999999999
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -7084.332999999999
$ws.Range("N122").Value = -10400.0002
$ws.Range("H132").Value = 5407834.5
$ws.Range("I132").Value = 6898960
$ws.Range("J132").Value = 2504.75
$ws.Range("K132").Value = 20696880
$ws.Range("L132").Value = 7514.25
$ws.Range("M132").Value = -20694350
$ws.Range("N132").Value = -12574.25
$ws.Range("H138").Value = 2041.0869
$ws.Range("I138").Value = 1310.238
$ws.Range("J138").Value = 3177.963
$ws.Range("K138").Value = 3930.714
$ws.Range("L138").Value = 9533.889000000001
$ws.Range("M138").Value = 1209.286
$ws.Range("N138").Value = -19813.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7104.8247
$ws.Range("I32").Value = 6134.1724
$ws.Range("J32").Value = 15549.5
$ws.Range("K32").Value = 6134.1724
$ws.Range("L32").Value = 15549.5
$ws.Range("M32").Value = -5847.1724
$ws.Range("N32").Value = -16123.5
$ws.Range("H45").Value = 1214.1072
$ws.Range("I45").Value = 1088
$ws.Range("J45").Value = 2265
$ws.Range("K45").Value = 1088
$ws.Range("L45").Value = 2265
$ws.Range("M45").Value = -711
$ws.Range("N45").Value = -3019
$ws.Range("H61").Value = 7588.8423
$ws.Range("I61").Value = 11266.546
$ws.Range("J61").Value = 2532
$ws.Range("K61").Value = 11266.546
$ws.Range("L61").Value = 2532
$ws.Range("M61").Value = -11054.546
$ws.Range("N61").Value = -2956
$ws.Range("H80").Value = 25405.285
$ws.Range("J80").Value = 25405.285
$ws.Range("L80").Value = 25405.285
$ws.Range("N80").Value = -27401.285
$ws.Range("H83").Value = 25405.285
$ws.Range("J83").Value = 25405.285
$ws.Range("L83").Value = 76215.855
$ws.Range("N83").Value = -86199.855
$ws.Range("H97").Value = 923.41174
$ws.Range("I97").Value = 856.125
$ws.Range("K97").Value = 856.125
$ws.Range("M97").Value = -360.125
$ws.Range("H102").Value = 2954.2856
$ws.Range("I102").Value = 2946.6667
$ws.Range("K102").Value = 2946.6667
$ws.Range("M102").Value = -1324.6667
$ws.Range("H122").Value = 2316.1177
$ws.Range("I122").Value = 1984.25
$ws.Range("J122").Value = 2611.111
$ws.Range("K122").Value = 5952.75
$ws.Range("L122").Value = 7833.333
$ws.Range("M122").Value = -3502.75
$ws.Range("N122").Value = -12733.333
$ws.Range("H132").Value = 1920.4062
$ws.Range("I132").Value = 1776.4706
$ws.Range("J132").Value = 2485.077
$ws.Range("K132").Value = 5329.4118
$ws.Range("L132").Value = 7455.231000000001
$ws.Range("M132").Value = -2799.4118
$ws.Range("N132").Value = -12515.231
$ws.Range("H136").Value = 7588.8423
$ws.Range("I136").Value = 11266.546
$ws.Range("J136").Value = 2532
$ws.Range("K136").Value = 33799.638
$ws.Range("L136").Value = 7596
$ws.Range("M136").Value = -31249.638
$ws.Range("N136").Value = -12696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 25888.889
$ws.Range("I126").Value = 26000
$ws.Range("J126").Value = 25875
$ws.Range("K126").Value = 26000
$ws.Range("L126").Value = 25875
$ws.Range("M126").Value = -21060
$ws.Range("N126").Value = -35755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1799.6
$ws.Range("I16").Value = 1642.1428
$ws.Range("J16").Value = 1937.375
$ws.Range("K16").Value = 1642.1428
$ws.Range("L16").Value = 1937.375
$ws.Range("M16").Value = -1355.1428
$ws.Range("N16").Value = -2511.375
$ws.Range("H58").Value = 17246332
$ws.Range("I58").Value = 2546.6875
$ws.Range("K58").Value = 2546.6875
$ws.Range("M58").Value = -2343.6875
$ws.Range("H99").Value = 1779.4667
$ws.Range("I99").Value = 1767.4286
$ws.Range("J99").Value = 1790
$ws.Range("K99").Value = 1767.4286
$ws.Range("L99").Value = 1790
$ws.Range("M99").Value = -269.4286
$ws.Range("N99").Value = -4786
$ws.Range("H105").Value = 1173.95
$ws.Range("I105").Value = 1498.3334
$ws.Range("J105").Value = 687.375
$ws.Range("K105").Value = 1498.3334
$ws.Range("L105").Value = 687.375
$ws.Range("M105").Value = 248.6666
$ws.Range("N105").Value = -4181.375
$ws.Range("H113").Value = 1799.6
$ws.Range("I113").Value = 1642.1428
$ws.Range("J113").Value = 1937.375
$ws.Range("K113").Value = 1642.1428
$ws.Range("L113").Value = 1937.375
$ws.Range("M113").Value = 527.8571999999999
$ws.Range("N113").Value = -6277.375
$ws.Range("H122").Value = 2516.0476
$ws.Range("J122").Value = 2026.25
$ws.Range("L122").Value = 6078.75
$ws.Range("N122").Value = -10978.75
$ws.Range("H126").Value = 1779.4667
$ws.Range("I126").Value = 1767.4286
$ws.Range("J126").Value = 1790
$ws.Range("K126").Value = 5302.2858
$ws.Range("L126").Value = 5370
$ws.Range("M126").Value = -2832.2858
$ws.Range("N126").Value = -10310
$ws.Range("H134").Value = 4600.56
$ws.Range("I134").Value = 4895.905
$ws.Range("J134").Value = 3050
$ws.Range("K134").Value = 14687.715
$ws.Range("L134").Value = 9150
$ws.Range("M134").Value = -12152.715
$ws.Range("N134").Value = -14220
$ws.Range("H136").Value = 17246332
$ws.Range("I136").Value = 2546.6875
$ws.Range("K136").Value = 7640.0625
$ws.Range("M136").Value = -5090.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1173.0934
$ws.Range("J131").Value = 924.3143
$ws.Range("L131").Value = 2772.9429
$ws.Range("N131").Value = -12852.9429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1591.4348
$ws.Range("I113").Value = 1782.7778
$ws.Range("J113").Value = 902.6
$ws.Range("K113").Value = 1782.7778
$ws.Range("L113").Value = 902.6
$ws.Range("M113").Value = 387.2221999999999
$ws.Range("N113").Value = -5242.6
$ws.Range("H132").Value = 3385.6365
$ws.Range("I132").Value = 3623.7917
$ws.Range("J132").Value = 3099.85
$ws.Range("K132").Value = 10871.3751
$ws.Range("L132").Value = 9299.549999999999
$ws.Range("M132").Value = -8341.375100000001
$ws.Range("N132").Value = -14359.55

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1583.9166
$ws.Range("I7").Value = 1285.8462
$ws.Range("K7").Value = 1285.8462
$ws.Range("M7").Value = -1173.8462
$ws.Range("H61").Value = 1927.1428
$ws.Range("I61").Value = 1898.3334
$ws.Range("J61").Value = 2100
$ws.Range("K61").Value = 1898.3334
$ws.Range("L61").Value = 2100
$ws.Range("M61").Value = -1696.3334
$ws.Range("N61").Value = -2504
$ws.Range("H113").Value = 1927.1428
$ws.Range("I113").Value = 1898.3334
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1898.3334
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 271.6666
$ws.Range("N113").Value = -6440
$ws.Range("H126").Value = 1583.9166
$ws.Range("I126").Value = 1285.8462
$ws.Range("K126").Value = 3857.5386
$ws.Range("M126").Value = -1387.5386
$ws.Range("H132").Value = 5281.4614
$ws.Range("I132").Value = 2098.6
$ws.Range("J132").Value = 8631.842000000001
$ws.Range("K132").Value = 6295.799999999999
$ws.Range("L132").Value = 25895.526
$ws.Range("M132").Value = -3765.799999999999
$ws.Range("N132").Value = -30955.526
$ws.Range("H136").Value = 2564.8333
$ws.Range("I136").Value = 2305.5833
$ws.Range("K136").Value = 6916.749899999999
$ws.Range("M136").Value = -4366.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1096.3043
$ws.Range("I113").Value = 497.46667
$ws.Range("J113").Value = 2219.125
$ws.Range("K113").Value = 1492.40001
$ws.Range("L113").Value = 6657.375
$ws.Range("M113").Value = 677.5999899999999
$ws.Range("N113").Value = -10997.375
$ws.Range("H122").Value = 3110.087
$ws.Range("I122").Value = 2747.4
$ws.Range("J122").Value = 3389.077
$ws.Range("K122").Value = 8242.200000000001
$ws.Range("L122").Value = 10167.231
$ws.Range("M122").Value = -5792.200000000001
$ws.Range("N122").Value = -15067.231
$ws.Range("H126").Value = 1919.3243
$ws.Range("I126").Value = 1655
$ws.Range("J126").Value = 2544.0908
$ws.Range("K126").Value = 4965
$ws.Range("L126").Value = 7632.2724
$ws.Range("M126").Value = -2495
$ws.Range("N126").Value = -12572.2724
$ws.Range("H132").Value = 14894.35
$ws.Range("I132").Value = 2245.2144
$ws.Range("J132").Value = 44409
$ws.Range("K132").Value = 6735.6432
$ws.Range("L132").Value = 133227
$ws.Range("M132").Value = -4205.6432
$ws.Range("N132").Value = -138287
